# DOMA-1155 regenerate excel templates
#
# Old layout (A1:H3):
#   A Дата снятия | B Адрес | C Услуга | D Номер прибора | E Место |
#   F Показание | G Житель | H Источник
#
# New layout (A1:L3):
#   A Дата снятия | B Адрес | C Квартира | D Услуга | E Номер прибора |
#   F Место | G..J Показание по тарифу №1..№4 | K Житель | L Источник

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Квартира" column right after "Адрес" (old C..H shift to D..I).
$ws.Columns("C").Insert()

# The old "Показание" column is now F. Split it into 4 tariff columns by
# inserting 3 more columns right after it (old G,H "Житель"/"Источник"
# end up shifting from G,H to K,L).
$ws.Columns("G").Insert()
$ws.Columns("G").Insert()
$ws.Columns("G").Insert()

# --- Row 1: headers ---
$ws.Range("A1").Value = "Дата снятия"
$ws.Range("B1").Value = "Адрес"
$ws.Range("C1").Value = "Квартира"
$ws.Range("D1").Value = "Услуга"
$ws.Range("E1").Value = "Номер прибора"
$ws.Range("F1").Value = "Место"
$ws.Range("G1").Value = "Показание по тарифу №1"
$ws.Range("H1").Value = "Показание по тарифу №2"
$ws.Range("I1").Value = "Показание по тарифу №3"
$ws.Range("J1").Value = "Показание по тарифу №4"
$ws.Range("K1").Value = "Житель"
$ws.Range("L1").Value = "Источник"

# --- Row 2: template placeholders for d.meter[i] ---
$ws.Range("A2").Value = "{d.meter[i].date}"
$ws.Range("B2").Value = "{d.meter[i].address}"
$ws.Range("C2").Value = "{d.meter[i].unitName}"
$ws.Range("D2").Value = "{d.meter[i].resource}"
$ws.Range("E2").Value = "{d.meter[i].number}"
$ws.Range("F2").Value = "{d.meter[i].place}"
$ws.Range("G2").Value = "{d.meter[i].value1}"
$ws.Range("H2").Value = "{d.meter[i].value2}"
$ws.Range("I2").Value = "{d.meter[i].value3}"
$ws.Range("J2").Value = "{d.meter[i].value4}"
$ws.Range("K2").Value = "{d.meter[i].clientName}"
$ws.Range("L2").Value = "{d.meter[i].source}"

# --- Row 3: template placeholders for d.meter[i + 1] ---
$ws.Range("A3").Value = "{d.meter[i + 1].date}"
$ws.Range("B3").Value = "{d.meter[i + 1].address}"
$ws.Range("C3").Value = "{d.meter[i + 1].unitName}"
$ws.Range("D3").Value = "{d.meter[i + 1].resource}"
$ws.Range("E3").Value = "{d.meter[i + 1].number}"
$ws.Range("F3").Value = "{d.meter[i + 1].place}"
$ws.Range("G3").Value = "{d.meter[i + 1].value1}"
$ws.Range("H3").Value = "{d.meter[i + 1].value2}"
$ws.Range("I3").Value = "{d.meter[i + 1].value3}"
$ws.Range("J3").Value = "{d.meter[i + 1].value4}"
$ws.Range("K3").Value = "{d.meter[i + 1].clientName}"
$ws.Range("L3").Value = "{d.meter[i + 1].source}"
